$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the report title (table caption) and the "Rolling 12 Months"
#    caption from "October" to "November" (EIA rolled the report forward one
#    month: 2006 - October 2016  ->  2006 - November 2016).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Table 4.1. Receipts, Average Cost, and Quality of Fossil Fuels: Total (All Sectors), 2006 - November 2016"

# ---------------------------------------------------------------------------
# 2) Insert a new row for "November" month data right after the existing
#    "October" row (row 52) and before the "Year to Date" summary section
#    (old row 53). This pushes everything from old row 53 onward down by one.
# ---------------------------------------------------------------------------
$ws.Rows("53:53").Insert()

# Copy the formatting of the October data row (now row 52) onto the new,
# blank November row (53) so the new row's styles match the rest of the
# monthly-data rows instead of Excel's default insert formatting.
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the November monthly figures.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 1082182
$ws.Range("C53").Value = 56396
$ws.Range("D53").Value = 2.09
$ws.Range("E53").Value = 40.03
$ws.Range("F53").Value = 1.29
$ws.Range("G53").Value = 114.7
$ws.Range("H53").Value = 9430
$ws.Range("I53").Value = 1551
$ws.Range("J53").Value = 10.07
$ws.Range("K53").Value = 61.25
$ws.Range("L53").Value = 0.49
$ws.Range("M53").Value = 88.9

# ---------------------------------------------------------------------------
# 3) Refresh the "Year to Date" annual totals (rows 55-57, formerly 54-56
#    before the insert) with the newly-published figures that now include
#    November.
# ---------------------------------------------------------------------------
# Year 2014
$ws.Range("A55").Value = 2014
$ws.Range("B55").Value = 15141009
$ws.Range("C55").Value = 779135
$ws.Range("D55").Value = 2.35
$ws.Range("E55").Value = 45.73
$ws.Range("F55").Value = 1.32
$ws.Range("G55").Value = 97.1
$ws.Range("H55").Value = 154066
$ws.Range("I55").Value = 25537
$ws.Range("J55").Value = 20.6
$ws.Range("K55").Value = 124.39
$ws.Range("L55").Value = 0.46
$ws.Range("M55").Value = 78.1

# Year 2015
$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 13961955
$ws.Range("C56").Value = 724360
$ws.Range("D56").Value = 2.23
$ws.Range("E56").Value = 42.97
$ws.Range("F56").Value = 1.29
$ws.Range("G56").Value = 102.8
$ws.Range("H56").Value = 137610
$ws.Range("I56").Value = 22664
$ws.Range("J56").Value = 11.68
$ws.Range("K56").Value = 70.98
$ws.Range("L56").Value = 0.49
$ws.Range("M56").Value = 74.9

# Year 2016 (preliminary, through November)
$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 11242072
$ws.Range("C57").Value = 581589
$ws.Range("D57").Value = 2.12
$ws.Range("E57").Value = 40.97
$ws.Range("F57").Value = 1.35
$ws.Range("G57").Value = 92.9
$ws.Range("H57").Value = 92662
$ws.Range("I57").Value = 15286
$ws.Range("J57").Value = 9.24
$ws.Range("K57").Value = 56.06
$ws.Range("L57").Value = 0.48
$ws.Range("M57").Value = 70.7

# ---------------------------------------------------------------------------
# 4) Update the "Rolling 12 Months Ending in October" section -> "November",
#    and refresh its rolling-12-month totals (rows 59-60, formerly 58-59).
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Rolling 12 months ending November 2015
$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 15415668
$ws.Range("C59").Value = 799784
$ws.Range("D59").Value = 2.26
$ws.Range("E59").Value = 43.48
$ws.Range("F59").Value = 1.29
$ws.Range("G59").Value = 103.3
$ws.Range("H59").Value = 155966
$ws.Range("I59").Value = 25641
$ws.Range("J59").Value = 11.94
$ws.Range("K59").Value = 72.7
$ws.Range("L59").Value = 0.49
$ws.Range("M59").Value = 79.6

# Rolling 12 months ending November 2016
$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 12366326
$ws.Range("C60").Value = 640158
$ws.Range("D60").Value = 2.12
$ws.Range("E60").Value = 41.01
$ws.Range("F60").Value = 1.34
$ws.Range("G60").Value = 94.5
$ws.Range("H60").Value = 102699
$ws.Range("I60").Value = 16943
$ws.Range("J60").Value = 9.2
$ws.Range("K60").Value = 55.81
$ws.Range("L60").Value = 0.47
$ws.Range("M60").Value = 72.4
